$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 23:52"

# Update the Cataluña row (row 5) statistics
$ws.Range("B5").Value = 32984
$ws.Range("C5").Value = 14298
$ws.Range("D5").Value = 15355
$ws.Range("E5").Value = 3331
